$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append a new test-mail row (row 7) ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A7").Value = "Wil je deze klant bellen?"
$ws.Range("B7").Value = "mailmind.test@zohomail.eu"
$ws.Range("C7").Value = "Testmail #5: Wil je deze klant bellen?"
$ws.Range("D7").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E7").Value = "Geachte afzender,`nDank voor uw e-mail. Om u beter van dienst te kunnen zijn, zouden we graag wat meer details ontvangen over welke klant we dienen te benaderen en waarvoor. Kunt u ons de naam van de klant en de reden voor het contact geven? Op die manier kunnen we dit efficiënt afhandelen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("F7").Value = "2025-07-29 21:37:31"
$ws.Range("G7").Value = "Ja"
$ws.Range("H7").Value = "Nee"
$ws.Range("I7").Value = "Ja"
$ws.Range("J7").Value = "Nee"

# Setting a multi-line value auto-expands the row height; reset it back to
# the sheet's default (un-"customized") height to match the source row.
$ws.Rows.Item(7).EntireRow.AutoFit()

# --- Extend conditional formatting ranges from row 6 to row 7 ---
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "6")
    $newRange = $ws.Range($col + "2:" + $col + "7")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- "Dashboard" sheet: update the category count for the new row's category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 2
